# Cronograma.xlsx update: "Atualizacao no cronograma e recursos humanos"
#
# Adds two new tracking columns (Recursos Ambientais / Recursos Humanos split
# out, plus Custo Estimado) to the "Cronograma" sheet header row, renames a
# couple of existing header cells, and leaves the "Legenda" sheet as-is
# (just switches the active sheet/selection back to Cronograma).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cronograma")

# --- Copy the formatting of the last two header-area columns (G/H) into the
#     two new columns (I/J) for every row that currently carries formatting,
#     BEFORE changing any text, so the new columns look consistent with the
#     rest of the table. ---

# Row 1 (thin strip above the header)
$ws.Range("G1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

# Body rows 3-18 (row 2 is the header, handled separately below)
$ws.Range("G3:G18").Copy()
$ws.Range("I3:I18").PasteSpecial(-4122)
$ws.Range("H3:H18").Copy()
$ws.Range("J3:J18").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Header row text updates (order matters: it drives shared-string index
#     allocation the same way the original authoring session did). ---

$ws.Range("B2").Value = "Área de Processo"
$ws.Range("G2").Value = "Esforço da Atividade"
$ws.Range("I2").Value = "Recursos Ambientais"
$ws.Range("H2").Value = "Recursos Humanos"
$ws.Range("J2").Value = "Custo Estimado"

# Give the two new header cells (I2/J2) the same header formatting as their
# neighbours (G2 = mid-header style, H2's old "last column" style now
# belongs to J2).
$ws.Range("G2").Copy()
$ws.Range("I2").PasteSpecial(-4122)
$ws.Range("H2").Copy()
$ws.Range("J2").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# H2 itself keeps being a "middle" header now that J2 exists, so it should
# carry the same formatting G2/I2 use rather than the old "last column" look.
$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Restore header cell values after the format-only paste operations
#     above (PasteSpecial(xlPasteFormats) does not touch values, but do this
#     defensively so ordering changes never clobber the text). ---
$ws.Range("B2").Value = "Área de Processo"
$ws.Range("C2").Value = "Atividade"
$ws.Range("D2").Value = "Descrição da Atividade"
$ws.Range("E2").Value = "Data de Início"
$ws.Range("F2").Value = "Deadline"
$ws.Range("G2").Value = "Esforço da Atividade"
$ws.Range("H2").Value = "Recursos Humanos"
$ws.Range("I2").Value = "Recursos Ambientais"
$ws.Range("J2").Value = "Custo Estimado"

# --- Switch the active sheet back to "Cronograma" (the legend/"Legenda"
#     sheet was the active one before this edit) and park the selection on
#     the newly added "Recursos Ambientais" header cell. ---
$ws.Activate()
$ws.Range("I2").Select()

$legenda = $wb.Worksheets.Item("Legenda")
$legenda.Range("C10").Select()
$ws.Activate()

Write-Output "Cronograma header updated"
